# Adding patchmethod use case
# Add a new "UpdateUser" worksheet after the existing "AddUser" sheet and
# populate it with a small request payload used by the PATCH endpoint.

$wb = $excel.ActiveWorkbook

$addUserSheet = $wb.Worksheets.Item("AddUser")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$updateUserSheet = $wb.Worksheets.Add($null, $lastSheet)
$updateUserSheet.Name = "UpdateUser"

# Header row
$updateUserSheet.Range("A1").Value = "name"
$updateUserSheet.Range("B1").Value = "email"
$updateUserSheet.Range("C1").Value = "userId"

# Data row
$updateUserSheet.Range("A2").Value = "test"
$updateUserSheet.Range("B2").Value = "test"
$updateUserSheet.Range("C2").Value = "'7354294"

# Move the view's selection off the AddUser sheet, onto the freshly added one
$addUserSheet.Range("A1:D2").Select() | Out-Null
$updateUserSheet.Range("E17").Select() | Out-Null
$updateUserSheet.Activate() | Out-Null
